# Fixing the big mistake: correct the Total (B) and Community (D)
# monthly consumption values for months 1-12 (rows 2-13) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 2;  B = 11766.92527505004;  D = 792.6157416166667 },
    @{ Row = 3;  B = 10948.09084975003;  D = 732.1203122666667 },
    @{ Row = 4;  B = 11727.08913428337;  D = 776.3483170833333 },
    @{ Row = 5;  B = 11334.08435408337;  D = 758.5950234166667 },
    @{ Row = 6;  B = 11771.38486485004;  D = 778.4913420666666 },
    @{ Row = 7;  B = 11347.8633005167;   D = 773.5021007 },
    @{ Row = 8;  B = 11737.61626008337;  D = 768.0996165 },
    @{ Row = 9;  B = 11741.40859793337;  D = 787.38040185 },
    @{ Row = 10; B = 11349.54056110004;  D = 749.86498935 },
    @{ Row = 11; B = 11743.28610040004;  D = 784.8573524333333 },
    @{ Row = 12; B = 11369.75652638337;  D = 757.8676801666667 },
    @{ Row = 13; B = 11315.1905300167;   D = 759.4010134499999 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
